# Applies updated vm_pu values for case with 380 kV (Case_2_171)
# Updates rows 2-25 (bus indices 0-23), columns B,C,D,E,F,I,J,K,L,M,N

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")
$colIdx = @{"B"=2;"C"=3;"D"=4;"E"=5;"F"=6;"I"=9;"J"=10;"K"=11;"L"=12;"M"=13;"N"=14}

$data = @(
    @(1.02,1.044563904612918,1.050569963871977,1.056615208865503,1.063328348071551,1.042559004350166,1.049628316352711,1.053323874959923,1.059352434032566,1.066047284487678,1.020393815394375),
    @(1.02,1.045500130253592,1.051303723881312,1.05756075773059,1.064261744130012,1.042789179414156,1.050211914009702,1.053870310581533,1.060111321250791,1.066795372704111,1.020590352569498),
    @(1.02,1.046106301772221,1.051778886939308,1.058173828247076,1.064866690860507,1.042937135083555,1.050589274236392,1.05422359261112,1.060602981310505,1.067279781652605,1.02071736895573),
    @(1.02,1.046361224082483,1.051978733027154,1.058431857876397,1.065121242973447,1.042999099787995,1.05074785173712,1.054372039972501,1.060809819881572,1.067483508779871,1.020770728837851),
    @(1.02,1.04640403178009,1.0520122931589,1.058475199432259,1.065163996986876,1.043009490088299,1.050774473773618,1.054396960663526,1.060844557437701,1.067517720211765,1.020779685968633),
    @(1.02,1.046109707712833,1.051781556950232,1.058177274895058,1.064870091286459,1.042937963986389,1.050591393412271,1.054225576457745,1.060605744530591,1.067282503545286,1.020718082102124),
    @(1.02,1.04488022904856,1.050817863821349,1.056934504525384,1.063643590486661,1.042636996305927,1.049825600522896,1.053508606625729,1.05960877665347,1.066300032075444,1.020460268185325),
    @(1.02,1.042716615432952,1.049122616201734,1.05475412143734,1.061489884611289,1.042099147720163,1.048474176306849,1.052242972270278,1.057856710526505,1.064571493291316,1.020004787031799),
    @(1.02,1.041276204106728,1.047994482265507,1.053307027648844,1.060059236706445,1.041735567279536,1.047571938062598,1.051397766320088,1.056691908905542,1.063421020828241,1.019700362262414),
    @(1.02,1.040652975785024,1.047506486160845,1.052681977832701,1.059440990596968,1.041576949983155,1.047180964730293,1.051031450921385,1.056188320576029,1.062923315701244,1.01956836491176),
    @(1.02,1.040421553945389,1.04732529811328,1.052450040946865,1.059211533211356,1.041517854981578,1.047035695911763,1.050895335578932,1.056001383555489,1.062738515675676,1.019519308632694),
    @(1.02,1.0404711914018,1.047364160129876,1.052499781564926,1.059260744126751,1.041530539083216,1.047066858539402,1.050924535006603,1.056041476817831,1.062778152730477,1.019529832568211),
    @(1.02,1.04063384489103,1.047491507553875,1.052662801066774,1.059422019760375,1.041572068789473,1.047168957655133,1.051020200588955,1.056172865887567,1.062908038642516,1.019564310444635),
    @(1.02,1.040734070762085,1.047569980539865,1.052763273869745,1.059521411787852,1.041597633109734,1.047231858497532,1.051079136781241,1.056253834714707,1.062988074925936,1.019585549890386),
    @(1.02,1.041317575892687,1.048026879466093,1.05334854294887,1.060100293777883,1.041746069271624,1.047597879467741,1.05142207047467,1.05672534690401,1.063454061606118,1.019709118743284),
    @(1.02,1.041683722198276,1.048313613099822,1.053716083087756,1.060463742466731,1.041838862559116,1.047827395537573,1.051637094517254,1.057021323271671,1.06374648591597,1.019786582445908),
    @(1.02,1.04189733532726,1.048480907488575,1.053930612622336,1.060675854896743,1.04189287298713,1.047961239544468,1.051762481955806,1.057194036311963,1.063917096088853,1.019831748377179),
    @(1.02,1.041970179690206,1.048537958570586,1.054003786980948,1.06074819988971,1.041911269740604,1.048006871974296,1.051805230300003,1.057252939687191,1.063975277180404,1.019847145840845),
    @(1.02,1.041644433381329,1.048282844370704,1.0536766340265,1.060424735547669,1.041828918539333,1.047802773611346,1.051614027821898,1.056989560039207,1.063715106993419,1.019778273114183),
    @(1.02,1.040585945472598,1.047454004818649,1.052614789403721,1.05937452293981,1.041559844221299,1.047138893215416,1.050992030818101,1.056134171811707,1.06286978854212,1.019554158300298),
    @(1.02,1.039920854387924,1.046933317583688,1.051948522660767,1.058715293730852,1.041389639946591,1.046721231625498,1.05060067072618,1.055597039187753,1.062338707804912,1.019413094787238),
    @(1.02,1.040273391426871,1.047209301868312,1.052301594168849,1.059064660669376,1.041479965594965,1.046942665680919,1.05080816490235,1.055881718142756,1.062620205022026,1.019487889663869),
    @(1.02,1.041662186163891,1.048296747288303,1.05369445889715,1.060442360727933,1.041833412167127,1.047813899289134,1.051624450768977,1.057003912245001,1.063729285641924,1.019782027797119),
    @(1.02,1.043275613570529,1.049560526162298,1.055316664019496,1.062045766509188,1.042239080648844,1.048823783081921,1.052570428835151,1.058309095237997,1.065018034615367,1.020122677330079)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $col = $colIdx[$cols[$j]]
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
